# Añadido de página de dashboard con estadisticas de los tickets
#
# Net content changes described by the diff:
#  - Sheet "Estado Tickets": ticket-state "Abierto" renamed to "Solicitud".
#  - Sheet "Categorías de Tickets":
#      * New header row inserted at the very top: A1 = "Categoria",
#        B1 = "Sub Categoria".
#      * The old top category cell (A2, "Incidencias / Errores") is
#        replaced with a placeholder value "_".
#      * The "Integraciones entre sistemas" sub-category row is removed,
#        shifting the remaining rows below it up by one.
#  - Selections: sheet1 ends up with B6 selected (and stays the active
#    tab), sheet2 ends up with C39 selected.

$wb  = $excel.ActiveWorkbook

$wsCategorias = $wb.Worksheets.Item("Categorías de Tickets")
$wsEstados    = $wb.Worksheets.Item("Estado Tickets")

# --- Estado Tickets: "Abierto" -> "Solicitud" ---------------------------
$wsEstados.Range("B2").Value = "Solicitud"

# --- Categorías de Tickets: new header row + placeholder + row removal --
$wsCategorias.Range("A1").Value = "Categoria"
$wsCategorias.Range("A2").Value = "_"
$wsCategorias.Range("B1").Value = "Sub Categoria"

# Remove the "Integraciones entre sistemas" row (row 22), shifting
# everything below it up by one row.
$wsCategorias.Rows.Item(22).Delete()

# The trailing row keeps the same value as the row above it (matches the
# source edit, which leaves a duplicated entry at the bottom of the list).
$wsCategorias.Range("B36").Value = $wsCategorias.Range("B35").Value2

# --- Selections -----------------------------------------------------------
$wsEstados.Range("C39").Select() | Out-Null

$wsCategorias.Activate() | Out-Null
$wsCategorias.Range("B6").Select() | Out-Null
